$d = $word.ActiveDocument

# The two list items describing finished functionalities need to be
# struck through, matching the formatting already used on all the other
# completed items in the list ("edita detaliile profilului" and
# "I se salveaza/salvează ultimele 10 filme ...").
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if (($t -like "*edita detaliile profilului*") -or ($t -like "*ultimele 10 filme*")) {
        $p.Range.Font.StrikeThrough = $true
    }
}
